$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before column U (shifts old V,W -> Z,AA)
$ws.Range("U1:X26").EntireColumn.Insert()

# Set the new header labels in row 3
$ws.Range("U3").Value = "Running time (seg)"
$ws.Range("V3").Value = "Data volume (GB)"
$ws.Range("W3").Value = "Sample size"
$ws.Range("X3").Value = "Ranking"

# Extend the header merge from R1:T2 to R1:X2
$ws.Range("R1:T2").UnMerge()
$ws.Range("R1:X2").Merge()
